$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1503.614180982113
$ws.Range("C2").Value = 0.2504771249345515
$ws.Range("D2").Value = 0.9841745793769782
$ws.Range("E2").Value = 0.915708021093533
